$d = $word.ActiveDocument

# Anchor an empty (collapsed) range at the very end of the document body,
# right after the last paragraph ("...левом нижнем перекрестке.").
$insertionPoint = $d.Content.End
$r = $d.Range($insertionPoint, $insertionPoint)

# Build the new paragraph as raw WordprocessingML and insert it via
# Range.InsertXML. This lets us express the "Python"/"pygame" runs with
# an explicit w:lang="en-US" on their run properties, matching the target
# OOXML exactly (something that isn't reliably reachable purely through
# Font.LanguageID on this host). The paragraph mirrors the formatting
# (Times New Roman, sz 28, line spacing 360 auto, justify both) used by
# the preceding paragraph in the document.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
'<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Программа была выполнена на </w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Python</w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">с помощью </w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>pygame</w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r>' +
'</w:p></w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
